# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row => new F value)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    3  = 14747
    4  = 18099
    6  = 92
    8  = 220
    9  = 28
    10 = 54
    16 = 184
    18 = 1365
    19 = 150
    20 = 80
    22 = 220
    23 = 7487
    29 = 5900
    33 = 152
    34 = 243
    35 = 5194
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (row => new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    3  = 14747
    4  = 18099
    6  = 92
    8  = 220
    9  = 28
    10 = 54
    16 = 184
    18 = 1365
    19 = 150
    20 = 80
    23 = 220
    24 = 7487
    31 = 5900
    35 = 152
    36 = 243
    37 = 5194
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
